$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts A:E -> B:F and
# carries over the existing header style (s="1") on row 1 for B1:F1.
$ws.Range("A1").EntireColumn.Insert()

# Match the header formatting used by the other header cells onto the
# new A1 header cell (reuses the same style record instead of creating
# a new one).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New "ID" column header + row labels for the data rows.
$ids = @("ID","Hb 2","Hb 3","S 24","S 28","Hb 107","Hb 66","Hb 69","Hb 95","Hb 99","Hb 92","Hb 40","Hb 41","S 11","Hb 57","S 21","S 22","S 3","S 4","S 5","Hb 74","Hb 79","Hb 32","S 15","S 16")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
